# Updated cryptos list (price + 1h volume change) refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "69.852.40"
$ws.Range("E2").Value = "  +0.65%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.571.57"
$ws.Range("E3").Value = "  +0.20%  "

$ws.Range("E4").Value = "  +0.19%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "574.03"
$ws.Range("E5").Value = "  -2.62%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "190.26"
$ws.Range("E6").Value = "  -1.06%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.630"
$ws.Range("E7").Value = "  -1.89%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.567.08"
$ws.Range("E8").Value = "  +0.06%  "

$ws.Range("E9").Value = "  +0.03%  "

$ws.Range("E10").Value = "  -3.52%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.658"
$ws.Range("E11").Value = "  -0.44%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "56.42"
$ws.Range("E12").Value = "  -3.17%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000298"
$ws.Range("E13").Value = "  +1.79%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "9.76"
$ws.Range("E14").Value = "  +0.89%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.145.54"
$ws.Range("E15").Value = "  +0.43%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "20.10"
$ws.Range("E16").Value = "  +4.44%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.577.34"
$ws.Range("E17").Value = "  +0.45%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "69.750.16"
$ws.Range("E18").Value = "  +0.64%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.45"
$ws.Range("E19").Value = "  +0.17%  "

$ws.Range("E20").Value = "  +1.06%  "

$ws.Range("E21").Value = "  -0.76%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "19.95"
$ws.Range("E22").Value = "  +16.80%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "471.67"
$ws.Range("E23").Value = "  -6.31%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "5.10"
$ws.Range("E24").Value = "  -8.04%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "4.33"
$ws.Range("E25").Value = "  -2.25%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "88.43"
$ws.Range("E26").Value = "  -2.99%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "3.05"
$ws.Range("E27").Value = "  +0.20%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "11.10"
$ws.Range("E28").Value = "  -0.74%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.26"
$ws.Range("E29").Value = "  -0.48%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.74"
$ws.Range("E30").Value = "  +3.75%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "32.02"
$ws.Range("E31").Value = "  -0.05%  "

$ws.Range("E32").Value = "  +5.26%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "12.08"
$ws.Range("E33").Value = "  -0.66%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "66.04"
$ws.Range("E34").Value = "  +1.00%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "588.57"
$ws.Range("E35").Value = "  -4.82%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "39.55"
$ws.Range("E36").Value = "  +4.62%  "

$ws.Range("B37").Value = "Dai"
$ws.Range("C37").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.999"
$ws.Range("E37").Value = "  -0.12%  "

$ws.Range("B38").Value = "PEPE"
$ws.Range("C38").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0₃0802"
$ws.Range("E38").Value = "  -3.43%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.399"
$ws.Range("E39").Value = "  +0.45%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.144"
$ws.Range("E40").Value = "  -2.73%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.94"
$ws.Range("E41").Value = "  +8.64%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.53"
$ws.Range("E42").Value = "  -3.01%  "

$ws.Range("B43").Value = "Maker"
$ws.Range("C43").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.223.54"
$ws.Range("E43").Value = "  -3.02%  "

$ws.Range("B44").Value = "ThetaToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "3.15"
$ws.Range("E44").Value = "  +2.54%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.13"
$ws.Range("E45").Value = "  +7.08%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0444"
$ws.Range("E46").Value = "  +0.90%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "9.54"
$ws.Range("E47").Value = "  +5.06%  "

$ws.Range("E48").Value = "  +2.58%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.136"
$ws.Range("E49").Value = "  -0.81%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.00"
$ws.Range("E50").Value = "  +0.43%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "3.15"
$ws.Range("E51").Value = "  -2.40%  "
